# New crime data collected — weekly update for cs-en-us-106pct.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# --- Row 15 (Rape) ---------------------------------------------------------
$ws.Range("F15").Value = 1
$ws.Range("M15").Value = -40
$ws.Range("N15").Value = -40

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 15.384615384615
$ws.Range("I16").Value = 64
$ws.Range("J16").Value = 64
$ws.Range("L16").Value = 45.454545454545
$ws.Range("M16").Value = -18.987341772151
$ws.Range("N16").Value = -78.378378378378

# --- Row 17 (Fel. Assault) ---------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 8
$ws.Range("I17").Value = 98
$ws.Range("J17").Value = 86
$ws.Range("K17").Value = 13.953488372093
$ws.Range("L17").Value = 20.987654320987
$ws.Range("M17").Value = 164.864864864865
$ws.Range("N17").Value = -7.547169811320

# --- Row 18 (Burglary) -------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = -28.571428571428
$ws.Range("M18").Value = -63.414634146341
$ws.Range("N18").Value = -92.084432717678

# --- Row 19 (Gr. Larceny) ----------------------------------------------------
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 14.285714285714
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = 1.960784313725
$ws.Range("I19").Value = 177
$ws.Range("J19").Value = 229
$ws.Range("K19").Value = -22.707423580786
$ws.Range("L19").Value = 30.147058823529
$ws.Range("M19").Value = 86.315789473684
$ws.Range("N19").Value = 1.142857142857

# --- Row 20 (G.L.A.) ----------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("I20").Value = 73
$ws.Range("J20").Value = 93
$ws.Range("K20").Value = -21.505376344086
$ws.Range("L20").Value = 65.909090909090
$ws.Range("M20").Value = -16.091954022988
$ws.Range("N20").Value = -93.209302325581

# --- Row 21 (TOTAL) -----------------------------------------------------------
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 23.076923076923
$ws.Range("F21").Value = 109
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = -6.837606837606
$ws.Range("I21").Value = 448
$ws.Range("J21").Value = 520
$ws.Range("K21").Value = -13.846153846153
$ws.Range("L21").Value = 31.764705882352
$ws.Range("M21").Value = 14.871794871794
$ws.Range("N21").Value = -78.178275694106

# --- Row 22 (Transit) ---------------------------------------------------------
$ws.Range("M22").Value = -37.5

# --- Row 24 (Petit Larceny) ---------------------------------------------------
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -41.176470588235
$ws.Range("F24").Value = 95
$ws.Range("G24").Value = 133
$ws.Range("H24").Value = -28.571428571428
$ws.Range("I24").Value = 395
$ws.Range("J24").Value = 427
$ws.Range("K24").Value = -7.494145199063
$ws.Range("L24").Value = 66.666666666666
$ws.Range("M24").Value = 108.994708994709

# --- Row 25 (Misd. Assault) ---------------------------------------------------
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 26.315789473684
$ws.Range("I25").Value = 158
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = 28.455284552845
$ws.Range("L25").Value = 31.666666666666
$ws.Range("M25").Value = 3.267973856209

# --- Row 26 (UCR Rape*) --------------------------------------------------------
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100

# --- Row 27 (Other Sex Crimes) --------------------------------------------------
# D27/E27 flip from text ("0"/"***.*") to real numbers, so pick up the number
# formats used by sibling numeric/percent cells in the same row before writing.
$ws.Range("D27").NumberFormat = $ws.Range("C27").NumberFormat
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = $ws.Range("H27").NumberFormat
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 5.263157894736
$ws.Range("L27").Value = 42.857142857142

# --- Row 28 (Shooting Vic.) ------------------------------------------------------
$ws.Range("M28").Value = 33.333333333333
$ws.Range("N28").Value = -78.947368421052

# --- Row 29 (Shooting Inc.) ------------------------------------------------------
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -89.473684210526
